$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing spaces from the three sector-description cells
# (RUV_OECD "space in sector 4-5 correction")
$ws.Range("A11").Value = "Manufacture of coke and refined petroleum products"
$ws.Range("A12").Value = "Manufacture of chemicals and chemical products"
$ws.Range("A27").Value = "Sewerage; waste collection, treatment and disposal activities; materials recovery; remediation activities and other waste management services"

# Update the active selection to match the saved view state
$ws.Range("A13").Select()
